# Actualizacion Datos Personales 4 nov
# Update the "Promedio", "Blancos" and "Por_Blan" columns (I, J, K) for
# Pesce Bautista Victor Manuel's groups 1AV (row 10) and 1CV (row 12)
# on the "1er Parcial" and "3er Parcial" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("1er Parcial", "3er Parcial")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # Row 10 -> 1AV
    $ws.Range("I10").Value = 7.5
    $ws.Range("J10").Value = 3
    $ws.Range("K10").Value = 7.5

    # Row 12 -> 1CV
    $ws.Range("I12").Value = 8.5
    $ws.Range("J12").Value = 4
    $ws.Range("K12").Value = 9.300000000000001
}
